$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 107,5
$data[0,0] = "Giorgian De Arrascaeta"
$data[0,1] = "Flamengo"
$data[0,2] = 6
$data[0,3] = 6
$data[0,4] = 3
$data[1,0] = "Pablo Vegetti"
$data[1,1] = "Vasco da Gama"
$data[1,2] = 6
$data[1,3] = 4
$data[1,4] = 1
$data[2,0] = "Pedro Raul"
$data[2,1] = "Ceará"
$data[2,2] = 7
$data[2,3] = 4
$data[2,4] = $null
$data[3,0] = "Yuri Alberto"
$data[3,1] = "Corinthians"
$data[3,2] = 7
$data[3,3] = 4
$data[3,4] = 1
$data[4,0] = "Kaio Jorge"
$data[4,1] = "Cruzeiro"
$data[4,2] = 7
$data[4,3] = 4
$data[4,4] = 1
$data[5,0] = "Reinaldo"
$data[5,1] = "Mirassol"
$data[5,2] = 6
$data[5,3] = 4
$data[5,4] = $null
$data[6,0] = "Pedro"
$data[6,1] = "Flamengo"
$data[6,2] = 5
$data[6,3] = 4
$data[6,4] = 1
$data[7,0] = "Memphis Depay"
$data[7,1] = "Corinthians"
$data[7,2] = 6
$data[7,3] = 3
$data[7,4] = 2
$data[8,0] = "Ferreira"
$data[8,1] = "São Paulo"
$data[8,2] = 7
$data[8,3] = 3
$data[8,4] = 2
$data[9,0] = "Gabriel Barbosa"
$data[9,1] = "Cruzeiro"
$data[9,2] = 7
$data[9,3] = 2
$data[9,4] = $null
$data[10,0] = "Rony"
$data[10,1] = "Atlético Mineiro"
$data[10,2] = 6
$data[10,3] = 2
$data[10,4] = $null
$data[11,0] = "Gabriel Taliari"
$data[11,1] = "Juventude"
$data[11,2] = 4
$data[11,3] = 2
$data[11,4] = $null
$data[12,0] = "Álvaro Barreal"
$data[12,1] = "Santos"
$data[12,2] = 5
$data[12,3] = 2
$data[12,4] = $null
$data[13,0] = "Marllon"
$data[13,1] = "Ceará"
$data[13,2] = 7
$data[13,3] = 2
$data[13,4] = $null
$data[14,0] = "Ignacio Laquintana"
$data[14,1] = "Bragantino"
$data[14,2] = 6
$data[14,3] = 2
$data[14,4] = $null
$data[15,0] = "Igor Jesus"
$data[15,1] = "Botafogo"
$data[15,2] = 7
$data[15,3] = 2
$data[15,4] = 1
$data[16,0] = "Alan Patrick"
$data[16,1] = "Internacional"
$data[16,2] = 6
$data[16,3] = 2
$data[16,4] = 4
$data[17,0] = "Rafael Borré"
$data[17,1] = "Internacional"
$data[17,2] = 6
$data[17,3] = 2
$data[17,4] = $null
$data[18,0] = "José Manuel López"
$data[18,1] = "Palmeiras"
$data[18,2] = 7
$data[18,3] = 2
$data[18,4] = $null
$data[19,0] = "Joaquín Piquerez"
$data[19,1] = "Palmeiras"
$data[19,2] = 7
$data[19,3] = 2
$data[19,4] = $null
$data[20,0] = "Eduardo Sasha"
$data[20,1] = "Bragantino"
$data[20,2] = 6
$data[20,3] = 2
$data[20,4] = 1
$data[21,0] = "Matheuzinho"
$data[21,1] = "Vitória"
$data[21,2] = 6
$data[21,3] = 2
$data[21,4] = $null
$data[22,0] = "Jefferson Savarino"
$data[22,1] = "Botafogo"
$data[22,2] = 6
$data[22,3] = 2
$data[22,4] = 1
$data[23,0] = "André Silva"
$data[23,1] = "São Paulo"
$data[23,2] = 7
$data[23,3] = 2
$data[23,4] = 1
$data[24,0] = "Daniel Borges"
$data[24,1] = "Mirassol"
$data[24,2] = 6
$data[24,3] = 2
$data[24,4] = $null
$data[25,0] = "Facundo Torres Pérez"
$data[25,1] = "Palmeiras"
$data[25,2] = 7
$data[25,3] = 2
$data[25,4] = 1
$data[26,0] = "Enmerson Batalla"
$data[26,1] = "Juventude"
$data[26,2] = 6
$data[26,3] = 2
$data[26,4] = 1
$data[27,0] = "Victor Gabriel"
$data[27,1] = "Internacional"
$data[27,2] = 3
$data[27,3] = 2
$data[27,4] = $null
$data[28,0] = "Dudu"
$data[28,1] = "Cruzeiro"
$data[28,2] = 5
$data[28,3] = 1
$data[28,4] = $null
$data[29,0] = "Lucas Ramon"
$data[29,1] = "Mirassol"
$data[29,2] = 6
$data[29,3] = 1
$data[29,4] = 1
$data[30,0] = "Matías Arezo"
$data[30,1] = "Grêmio"
$data[30,2] = 6
$data[30,3] = 1
$data[30,4] = 1
$data[31,0] = "Edenilson"
$data[31,1] = "Grêmio"
$data[31,2] = 7
$data[31,3] = 1
$data[31,4] = $null
$data[32,0] = "Juan Lucero"
$data[32,1] = "Fortaleza"
$data[32,2] = 7
$data[32,3] = 1
$data[32,4] = $null
$data[33,0] = "Tinga"
$data[33,1] = "Fortaleza"
$data[33,2] = 3
$data[33,3] = 1
$data[33,4] = $null
$data[34,0] = "Bruno Henrique"
$data[34,1] = "Internacional"
$data[34,2] = 6
$data[34,3] = 1
$data[34,4] = $null
$data[35,0] = "Léo Pereira"
$data[35,1] = "Flamengo"
$data[35,2] = 7
$data[35,3] = 1
$data[35,4] = $null
$data[36,0] = "Nuno Moreira"
$data[36,1] = "Vasco da Gama"
$data[36,2] = 7
$data[36,3] = 1
$data[36,4] = 1
$data[37,0] = "Gilberto"
$data[37,1] = "Bahia"
$data[37,2] = 6
$data[37,3] = 1
$data[37,4] = $null
$data[38,0] = "Héctor Hernández"
$data[38,1] = "Corinthians"
$data[38,2] = 7
$data[38,3] = 1
$data[38,4] = $null
$data[39,0] = "Eric Ramires"
$data[39,1] = "Bragantino"
$data[39,2] = 6
$data[39,3] = 1
$data[39,4] = $null
$data[40,0] = "Matheus Araújo"
$data[40,1] = "Ceará"
$data[40,2] = 7
$data[40,3] = 1
$data[40,4] = 1
$data[41,0] = "Mateo Ponte"
$data[41,1] = "Botafogo"
$data[41,2] = 7
$data[41,3] = 1
$data[41,4] = $null
$data[42,0] = "Lima"
$data[42,1] = "Fluminense"
$data[42,2] = 6
$data[42,3] = 1
$data[42,4] = $null
$data[43,0] = "Isidro Pitta"
$data[43,1] = "Bragantino"
$data[43,2] = 6
$data[43,3] = 1
$data[43,4] = $null
$data[44,0] = "Martinelli"
$data[44,1] = "Fluminense"
$data[44,2] = 7
$data[44,3] = 1
$data[44,4] = $null
$data[45,0] = "Leandro Martínez"
$data[45,1] = "Fortaleza"
$data[45,2] = 5
$data[45,3] = 1
$data[45,4] = $null
$data[46,0] = "Cristian Renato"
$data[46,1] = "Mirassol"
$data[46,2] = 5
$data[46,3] = 1
$data[46,4] = 1
$data[47,0] = "Enner Valencia"
$data[47,1] = "Internacional"
$data[47,2] = 7
$data[47,3] = 1
$data[47,4] = $null
$data[48,0] = "Wellington Rato"
$data[48,1] = "Vitória"
$data[48,2] = 4
$data[48,3] = 1
$data[48,4] = $null
$data[49,0] = "Bruno Henrique"
$data[49,1] = "Flamengo"
$data[49,2] = 7
$data[49,3] = 1
$data[49,4] = $null
$data[50,0] = "Chrystian Barletta"
$data[50,1] = "Sport"
$data[50,2] = 7
$data[50,3] = 1
$data[50,4] = $null
$data[51,0] = "Erick"
$data[51,1] = "Bahia"
$data[51,2] = 7
$data[51,3] = 1
$data[51,4] = $null
$data[52,0] = "Thaciano"
$data[52,1] = "Santos"
$data[52,2] = 6
$data[52,3] = 1
$data[52,4] = $null
$data[53,0] = "Diego Pituca"
$data[53,1] = "Santos"
$data[53,2] = 7
$data[53,3] = 1
$data[53,4] = $null
$data[54,0] = "Luciano Juba"
$data[54,1] = "Bahia"
$data[54,2] = 7
$data[54,3] = 1
$data[54,4] = 3
$data[55,0] = "Aylon"
$data[55,1] = "Ceará"
$data[55,2] = 5
$data[55,3] = 1
$data[55,4] = $null
$data[56,0] = "Luiz Mandaca"
$data[56,1] = "Juventude"
$data[56,2] = 6
$data[56,3] = 1
$data[56,4] = $null
$data[57,0] = "Matheus Babi"
$data[57,1] = "Juventude"
$data[57,2] = 4
$data[57,3] = 1
$data[57,4] = $null
$data[58,0] = "Emiliano Martínez"
$data[58,1] = "Palmeiras"
$data[58,2] = 7
$data[58,3] = 1
$data[58,4] = $null
$data[59,0] = "Rayan"
$data[59,1] = "Vasco da Gama"
$data[59,2] = 7
$data[59,3] = 1
$data[59,4] = 1
$data[60,0] = "Gabriel"
$data[60,1] = "Mirassol"
$data[60,2] = 4
$data[60,3] = 1
$data[60,4] = $null
$data[61,0] = "Erick"
$data[61,1] = "Bahia"
$data[61,2] = 7
$data[61,3] = 1
$data[61,4] = 1
$data[62,0] = "Samuel Xavier"
$data[62,1] = "Fluminense"
$data[62,2] = 7
$data[62,3] = 1
$data[62,4] = $null
$data[63,0] = "Lucas Halter"
$data[63,1] = "Vitória"
$data[63,2] = 7
$data[63,3] = 1
$data[63,4] = $null
$data[64,0] = "Fausto Vera"
$data[64,1] = "Atlético Mineiro"
$data[64,2] = 6
$data[64,3] = 1
$data[64,4] = $null
$data[65,0] = "Igor Gomes"
$data[65,1] = "Atlético Mineiro"
$data[65,2] = 5
$data[65,3] = 1
$data[65,4] = $null
$data[66,0] = "Martin Braithwaite"
$data[66,1] = "Grêmio"
$data[66,2] = 5
$data[66,3] = 1
$data[66,4] = $null
$data[67,0] = "Juninho Capixaba"
$data[67,1] = "Bragantino"
$data[67,2] = 6
$data[67,3] = 1
$data[67,4] = $null
$data[68,0] = "Renê"
$data[68,1] = "Fluminense"
$data[68,2] = 7
$data[68,3] = 1
$data[68,4] = $null
$data[69,0] = "Jhon Arias"
$data[69,1] = "Fluminense"
$data[69,2] = 7
$data[69,3] = 1
$data[69,4] = 3
$data[70,0] = "Erick Pulgar"
$data[70,1] = "Flamengo"
$data[70,2] = 7
$data[70,3] = 1
$data[70,4] = $null
$data[71,0] = "Gonzalo Plata"
$data[71,1] = "Flamengo"
$data[71,2] = 6
$data[71,3] = 1
$data[71,4] = $null
$data[72,0] = "Danilo"
$data[72,1] = "Flamengo"
$data[72,2] = 5
$data[72,3] = 1
$data[72,4] = $null
$data[73,0] = "Zé Ivaldo"
$data[73,1] = "Santos"
$data[73,2] = 7
$data[73,3] = 1
$data[73,4] = $null
$data[74,0] = "Janderson Costa"
$data[74,1] = "Vitória"
$data[74,2] = 7
$data[74,3] = 1
$data[74,4] = $null
$data[75,0] = "Diogo Barbosa"
$data[75,1] = "Fortaleza"
$data[75,2] = 6
$data[75,3] = 1
$data[75,4] = $null
$data[76,0] = "Lucas Romero"
$data[76,1] = "Cruzeiro"
$data[76,2] = 7
$data[76,3] = 1
$data[76,4] = $null
$data[77,0] = "Sergio Oliveira"
$data[77,1] = "Sport"
$data[77,2] = 4
$data[77,3] = 1
$data[77,4] = $null
$data[78,0] = "Ênio"
$data[78,1] = "Juventude"
$data[78,2] = 5
$data[78,3] = 1
$data[78,4] = 1
$data[79,0] = "Iury"
$data[79,1] = "Mirassol"
$data[79,2] = 6
$data[79,3] = 1
$data[79,4] = 1
$data[80,0] = "Tiquinho Soares"
$data[80,1] = "Santos"
$data[80,2] = 7
$data[80,3] = 1
$data[80,4] = 1
$data[81,0] = "Tomás Cuello"
$data[81,1] = "Atlético Mineiro"
$data[81,2] = 6
$data[81,3] = 1
$data[81,4] = $null
$data[82,0] = "Germán Cano"
$data[82,1] = "Fluminense"
$data[82,2] = 6
$data[82,3] = 1
$data[82,4] = $null
$data[83,0] = "Lucas Braga"
$data[83,1] = "Vitória"
$data[83,2] = 6
$data[83,3] = 1
$data[83,4] = $null
$data[84,0] = "Deyverson"
$data[84,1] = "Fortaleza"
$data[84,2] = 8
$data[84,3] = 1
$data[84,4] = $null
$data[85,0] = "Jhonatan"
$data[85,1] = "Bragantino"
$data[85,2] = 6
$data[85,3] = 1
$data[85,4] = 1
$data[86,0] = "Éverton Ribeiro"
$data[86,1] = "Bahia"
$data[86,2] = 6
$data[86,3] = 1
$data[86,4] = $null
$data[87,0] = "Edson Carioca"
$data[87,1] = "Mirassol"
$data[87,2] = 5
$data[87,3] = 1
$data[87,4] = 2
$data[88,0] = "Hulk"
$data[88,1] = "Atlético Mineiro"
$data[88,2] = 6
$data[88,3] = 1
$data[88,4] = 1
$data[89,0] = "Pedro Henrique"
$data[89,1] = "Ceará"
$data[89,2] = 3
$data[89,3] = 1
$data[89,4] = $null
$data[90,0] = "Ryan Rodrigues"
$data[90,1] = "São Paulo"
$data[90,2] = 7
$data[90,3] = 1
$data[90,4] = $null
$data[91,0] = "Vitinho"
$data[91,1] = "Botafogo"
$data[91,2] = 7
$data[91,3] = 1
$data[91,4] = $null
$data[92,0] = "Éverton"
$data[92,1] = "Flamengo"
$data[92,2] = 7
$data[92,3] = 1
$data[92,4] = 1
$data[93,0] = "Kayky"
$data[93,1] = "Bahia"
$data[93,2] = 7
$data[93,3] = 1
$data[93,4] = $null
$data[94,0] = "Christian"
$data[94,1] = "Cruzeiro"
$data[94,2] = 5
$data[94,3] = 1
$data[94,4] = 1
$data[95,0] = "Jemerson"
$data[95,1] = "Grêmio"
$data[95,2] = 7
$data[95,3] = 1
$data[95,4] = $null
$data[96,0] = "Carlinhos"
$data[96,1] = "Vitória"
$data[96,2] = 4
$data[96,3] = 1
$data[96,4] = $null
$data[97,0] = "Deivid Washington"
$data[97,1] = "Santos"
$data[97,2] = 7
$data[97,3] = 1
$data[97,4] = $null
$data[98,0] = "Pablo Felipe"
$data[98,1] = "Sport"
$data[98,2] = 4
$data[98,3] = 1
$data[98,4] = $null
$data[99,0] = "Kevin Serna"
$data[99,1] = "Fluminense"
$data[99,2] = 7
$data[99,3] = 1
$data[99,4] = $null
$data[100,0] = "Everaldo"
$data[100,1] = "Fluminense"
$data[100,2] = 7
$data[100,3] = 1
$data[100,4] = $null
$data[101,0] = "Braian Aguirre"
$data[101,1] = "Internacional"
$data[101,2] = 7
$data[101,3] = 1
$data[101,4] = $null
$data[102,0] = "Thiago Maia"
$data[102,1] = "Internacional"
$data[102,2] = 7
$data[102,3] = 1
$data[102,4] = $null
$data[103,0] = "Igor Coronado"
$data[103,1] = "Corinthians"
$data[103,2] = 4
$data[103,3] = 1
$data[103,4] = $null
$data[104,0] = "Cauly"
$data[104,1] = "Bahia"
$data[104,2] = 7
$data[104,3] = 1
$data[104,4] = 2
$data[105,0] = "Vítor Roque"
$data[105,1] = "Palmeiras"
$data[105,2] = 7
$data[105,3] = 1
$data[105,4] = $null
$data[106,0] = "Cristian Olivera"
$data[106,1] = "Grêmio"
$data[106,2] = 4
$data[106,3] = 1
$data[106,4] = $null
$ws.Range("A2:E108").Value = $data
